$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's cells (and the column's default format) switch to Text ("@"),
# which is what the header cells A1/B1 pick up too.
$ws.Cells.NumberFormat = "@"

# Add the new second row. Re-assert the Text format directly on A2 *before*
# writing the value, so "3" is kept as a text string (shared string) rather
# than being parsed into a numeric literal.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "3"

# Move the active selection to A2 (matches the diff's selection change)
$ws.Range("A2").Select()
